$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0.2690646666666667
$ws.Range("N2").Value = 0.807194
$ws.Range("O2").Value = 0.0885168346903475
$ws.Range("P2").Value = 0.09916786583441374
$ws.Range("Q2").Value = 0.04080383607644445
$ws.Range("R2").Value = 0.367234524688
$ws.Range("S2").Value = 0.0885168346903475
$ws.Range("T2").Value = 0.09916786583441374

$ws.Range("O3").Value = 0.4538005197112895
$ws.Range("P3").Value = 0.5084053130881303
$ws.Range("S3").Value = 0.4538005197112895
$ws.Range("T3").Value = 0.5084053130881303

$ws.Range("M4").Value = 0.130073
$ws.Range("N4").Value = 0.390219
$ws.Range("O4").Value = 0.04279138684880302
$ws.Range("P4").Value = 0.04794037794884388
$ws.Range("Q4").Value = 0.01972565716533333
$ws.Range("R4").Value = 0.177530914488
$ws.Range("S4").Value = 0.04279138684880302
$ws.Range("T4").Value = 0.04794037794884388

$ws.Range("M5").Value = 0.9794285
$ws.Range("N5").Value = 1.958857
$ws.Range("O5").Value = 0.3222121718899608
$ws.Range("P5").Value = 0.2406554907058306
$ws.Range("Q5").Value = 0.1485309849773334
$ws.Range("R5").Value = 0.891185909864
$ws.Range("S5").Value = 0.3222121718899608
$ws.Range("T5").Value = 0.2406554907058306

$ws.Range("M6").Value = 0.2817166666666667
$ws.Range("N6").Value = 0.84515
$ws.Range("O6").Value = 0.09267908685959904
$ws.Range("P6").Value = 0.1038309524227816
$ws.Range("Q6").Value = 0.04272252031111112
$ws.Range("R6").Value = 0.3845026828
$ws.Range("S6").Value = 0.09267908685959904
$ws.Range("T6").Value = 0.1038309524227816
